$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 1499.5  # ALC!H32: 2000 -> 1499.5
$ws.Cells.Item(32, 10).Value = 1499.5  # ALC!J32: 2000 -> 1499.5
$ws.Cells.Item(32, 12).Value = 1499.5  # ALC!L32: 2000 -> 1499.5
$ws.Cells.Item(32, 14).Value = -2151.5  # ALC!N32: -2652 -> -2151.5
$ws.Cells.Item(49, 8).Value = 793.4  # ALC!H49: 740.6667 -> 793.4
$ws.Cells.Item(49, 9).Value = 474.5  # ALC!I49: 475.33334 -> 474.5
$ws.Cells.Item(49, 11).Value = 1423.5  # ALC!K49: 1426.00002 -> 1423.5
$ws.Cells.Item(49, 13).Value = -1287.5  # ALC!M49: -1290.00002 -> -1287.5
$ws.Cells.Item(62, 8).Value = 3574.9333  # ALC!H62: 3608.2666 -> 3574.9333
$ws.Cells.Item(62, 9).Value = 3590.3333  # ALC!I62: 3676.625 -> 3590.3333
$ws.Cells.Item(62, 10).Value = 3551.8333  # ALC!J62: 3530.1428 -> 3551.8333
$ws.Cells.Item(62, 11).Value = 3590.3333  # ALC!K62: 3676.625 -> 3590.3333
$ws.Cells.Item(62, 12).Value = 3551.8333  # ALC!L62: 3530.1428 -> 3551.8333
$ws.Cells.Item(62, 13).Value = -2966.3333  # ALC!M62: -3052.625 -> -2966.3333
$ws.Cells.Item(62, 14).Value = -4799.8333  # ALC!N62: -4778.1428 -> -4799.8333
$ws.Cells.Item(65, 8).Value = 3574.9333  # ALC!H65: 3608.2666 -> 3574.9333
$ws.Cells.Item(65, 9).Value = 3590.3333  # ALC!I65: 3676.625 -> 3590.3333
$ws.Cells.Item(65, 10).Value = 3551.8333  # ALC!J65: 3530.1428 -> 3551.8333
$ws.Cells.Item(65, 11).Value = 17951.6665  # ALC!K65: 18383.125 -> 17951.6665
$ws.Cells.Item(65, 12).Value = 17759.1665  # ALC!L65: 17650.714 -> 17759.1665
$ws.Cells.Item(65, 13).Value = -14831.6665  # ALC!M65: -15263.125 -> -14831.6665
$ws.Cells.Item(65, 14).Value = -23999.1665  # ALC!N65: -23890.714 -> -23999.1665
$ws.Cells.Item(100, 8).Value = 2997.0833  # ALC!H100: 2996.818 -> 2997.0833
$ws.Cells.Item(100, 9).Value = 2952.5  # ALC!I100: 2905 -> 2952.5
$ws.Cells.Item(100, 11).Value = 2952.5  # ALC!K100: 2905 -> 2952.5
$ws.Cells.Item(100, 13).Value = -2411.5  # ALC!M100: -2364 -> -2411.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 6743.8086  # ARM!H61: 5882.2 -> 6743.8086
$ws.Cells.Item(61, 9).Value = 3420.2942  # ARM!I61: 2873.0698 -> 3420.2942
$ws.Cells.Item(61, 10).Value = 15436.077  # ARM!J61: 16664.916 -> 15436.077
$ws.Cells.Item(61, 11).Value = 3420.2942  # ARM!K61: 2873.0698 -> 3420.2942
$ws.Cells.Item(61, 12).Value = 15436.077  # ARM!L61: 16664.916 -> 15436.077
$ws.Cells.Item(61, 13).Value = -3208.2942  # ARM!M61: -2661.0698 -> -3208.2942
$ws.Cells.Item(61, 14).Value = -15860.077  # ARM!N61: -17088.916 -> -15860.077
$ws.Cells.Item(74, 8).Value = 3299.4736  # ARM!H74: 3091.3115 -> 3299.4736
$ws.Cells.Item(74, 9).Value = 1535.9025  # ARM!I74: 1410.4889 -> 1535.9025
$ws.Cells.Item(74, 11).Value = 1535.9025  # ARM!K74: 1410.4889 -> 1535.9025
$ws.Cells.Item(74, 13).Value = -661.9024999999999  # ARM!M74: -536.4889000000001 -> -661.9024999999999
$ws.Cells.Item(77, 8).Value = 3299.4736  # ARM!H77: 3091.3115 -> 3299.4736
$ws.Cells.Item(77, 9).Value = 1535.9025  # ARM!I77: 1410.4889 -> 1535.9025
$ws.Cells.Item(77, 11).Value = 7679.5125  # ARM!K77: 7052.444500000001 -> 7679.5125
$ws.Cells.Item(77, 13).Value = -3311.5125  # ARM!M77: -2684.444500000001 -> -3311.5125
$ws.Cells.Item(124, 8).Value = 34857  # ARM!H124: 33771.4 -> 34857
$ws.Cells.Item(124, 10).Value = 34857  # ARM!J124: 33771.4 -> 34857
$ws.Cells.Item(124, 12).Value = 34857  # ARM!L124: 33771.4 -> 34857
$ws.Cells.Item(124, 14).Value = -44677  # ARM!N124: -43591.4 -> -44677
$ws.Cells.Item(125, 8).Value = 67462  # ARM!H125: 67527 -> 67462
$ws.Cells.Item(125, 10).Value = 67462  # ARM!J125: 67527 -> 67462
$ws.Cells.Item(125, 12).Value = 67462  # ARM!L125: 67527 -> 67462
$ws.Cells.Item(125, 14).Value = -77302  # ARM!N125: -77367 -> -77302
$ws.Cells.Item(132, 8).Value = 4532.755  # ARM!H132: 3570.3235 -> 4532.755
$ws.Cells.Item(132, 9).Value = 1786.2858  # ARM!I132: 1240.0454 -> 1786.2858
$ws.Cells.Item(132, 10).Value = 7608.8  # ARM!J132: 7842.5 -> 7608.8
$ws.Cells.Item(132, 11).Value = 5358.857400000001  # ARM!K132: 3720.1362 -> 5358.857400000001
$ws.Cells.Item(132, 12).Value = 22826.4  # ARM!L132: 23527.5 -> 22826.4
$ws.Cells.Item(132, 13).Value = -2828.857400000001  # ARM!M132: -1190.1362 -> -2828.857400000001
$ws.Cells.Item(132, 14).Value = -27886.4  # ARM!N132: -28587.5 -> -27886.4
$ws.Cells.Item(136, 8).Value = 6743.8086  # ARM!H136: 5882.2 -> 6743.8086
$ws.Cells.Item(136, 9).Value = 3420.2942  # ARM!I136: 2873.0698 -> 3420.2942
$ws.Cells.Item(136, 10).Value = 15436.077  # ARM!J136: 16664.916 -> 15436.077
$ws.Cells.Item(136, 11).Value = 10260.8826  # ARM!K136: 8619.2094 -> 10260.8826
$ws.Cells.Item(136, 12).Value = 46308.231  # ARM!L136: 49994.74800000001 -> 46308.231
$ws.Cells.Item(136, 13).Value = -7710.882599999999  # ARM!M136: -6069.2094 -> -7710.882599999999
$ws.Cells.Item(136, 14).Value = -51408.231  # ARM!N136: -55094.74800000001 -> -51408.231

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(124, 8).Value = 79790  # BSM!H124: 79800 -> 79790
$ws.Cells.Item(124, 10).Value = 79790  # BSM!J124: 79800 -> 79790
$ws.Cells.Item(124, 12).Value = 79790  # BSM!L124: 79800 -> 79790
$ws.Cells.Item(124, 14).Value = -89610  # BSM!N124: -89620 -> -89610
$ws.Cells.Item(134, 8).Value = 19373.355  # BSM!H134: 21950.71 -> 19373.355
$ws.Cells.Item(134, 9).Value = 2073.3044  # BSM!I134: 2386.725 -> 2073.3044
$ws.Cells.Item(134, 10).Value = 80588.92  # BSM!J134: 87164 -> 80588.92
$ws.Cells.Item(134, 11).Value = 6219.9132  # BSM!K134: 7160.174999999999 -> 6219.9132
$ws.Cells.Item(134, 12).Value = 241766.76  # BSM!L134: 261492 -> 241766.76
$ws.Cells.Item(134, 13).Value = -3684.9132  # BSM!M134: -4625.174999999999 -> -3684.9132
$ws.Cells.Item(134, 14).Value = -246836.76  # BSM!N134: -266562 -> -246836.76

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3909.0425  # CRP!H31: 3870.0613 -> 3909.0425
$ws.Cells.Item(31, 9).Value = 1478.9546  # CRP!I31: 1420.6818 -> 1478.9546
$ws.Cells.Item(31, 10).Value = 6047.52  # CRP!J31: 5865.852 -> 6047.52
$ws.Cells.Item(31, 11).Value = 1478.9546  # CRP!K31: 1420.6818 -> 1478.9546
$ws.Cells.Item(31, 12).Value = 6047.52  # CRP!L31: 5865.852 -> 6047.52
$ws.Cells.Item(31, 13).Value = -1183.9546  # CRP!M31: -1125.6818 -> -1183.9546
$ws.Cells.Item(31, 14).Value = -6637.52  # CRP!N31: -6455.852 -> -6637.52
$ws.Cells.Item(34, 8).Value = 3909.0425  # CRP!H34: 3870.0613 -> 3909.0425
$ws.Cells.Item(34, 9).Value = 1478.9546  # CRP!I34: 1420.6818 -> 1478.9546
$ws.Cells.Item(34, 10).Value = 6047.52  # CRP!J34: 5865.852 -> 6047.52
$ws.Cells.Item(34, 11).Value = 1478.9546  # CRP!K34: 1420.6818 -> 1478.9546
$ws.Cells.Item(34, 12).Value = 6047.52  # CRP!L34: 5865.852 -> 6047.52
$ws.Cells.Item(34, 13).Value = -1276.9546  # CRP!M34: -1218.6818 -> -1276.9546
$ws.Cells.Item(34, 14).Value = -6451.52  # CRP!N34: -6269.852 -> -6451.52
$ws.Cells.Item(57, 8).Value = 2922330.5  # CRP!H57: 2658487.8 -> 2922330.5
$ws.Cells.Item(57, 10).Value = 3245700.5  # CRP!J57: 2923136.5 -> 3245700.5
$ws.Cells.Item(57, 12).Value = 3245700.5  # CRP!L57: 2923136.5 -> 3245700.5
$ws.Cells.Item(57, 14).Value = -3246820.5  # CRP!N57: -2924256.5 -> -3246820.5
$ws.Cells.Item(58, 8).Value = 1685500.9  # CRP!H58: 1654868.9 -> 1685500.9
$ws.Cells.Item(58, 9).Value = 2021626.8  # CRP!I58: 2067574 -> 2021626.8
$ws.Cells.Item(58, 10).Value = 4871.1113  # CRP!J58: 4048.5454 -> 4871.1113
$ws.Cells.Item(58, 11).Value = 2021626.8  # CRP!K58: 2067574 -> 2021626.8
$ws.Cells.Item(58, 12).Value = 4871.1113  # CRP!L58: 4048.5454 -> 4871.1113
$ws.Cells.Item(58, 13).Value = -2021423.8  # CRP!M58: -2067371 -> -2021423.8
$ws.Cells.Item(58, 14).Value = -5277.1113  # CRP!N58: -4454.5454 -> -5277.1113
$ws.Cells.Item(75, 8).Value = 45000  # CRP!H75: 40000 -> 45000
$ws.Cells.Item(75, 10).Value = 50000  # CRP!J75: 0 -> 50000
$ws.Cells.Item(75, 12).Value = 50000  # CRP!L75: 0 -> 50000
$ws.Cells.Item(75, 14).Value = -51996  # CRP!N75: None -> -51996
$ws.Cells.Item(76, 8).Value = 8596  # CRP!H76: 9102.857 -> 8596
$ws.Cells.Item(76, 9).Value = 8596  # CRP!I76: 9102.857 -> 8596
$ws.Cells.Item(76, 11).Value = 8596  # CRP!K76: 9102.857 -> 8596
$ws.Cells.Item(76, 13).Value = -8281  # CRP!M76: -8787.857 -> -8281
$ws.Cells.Item(78, 8).Value = 45000  # CRP!H78: 40000 -> 45000
$ws.Cells.Item(78, 10).Value = 50000  # CRP!J78: 0 -> 50000
$ws.Cells.Item(78, 12).Value = 150000  # CRP!L78: 0 -> 150000
$ws.Cells.Item(78, 14).Value = -159984  # CRP!N78: None -> -159984
$ws.Cells.Item(79, 8).Value = 8596  # CRP!H79: 9102.857 -> 8596
$ws.Cells.Item(79, 9).Value = 8596  # CRP!I79: 9102.857 -> 8596
$ws.Cells.Item(79, 11).Value = 8596  # CRP!K79: 9102.857 -> 8596
$ws.Cells.Item(79, 13).Value = -7504  # CRP!M79: -8010.857 -> -7504
$ws.Cells.Item(99, 8).Value = 2303.1667  # CRP!H99: 2302.6667 -> 2303.1667
$ws.Cells.Item(99, 9).Value = 2311.111  # CRP!I99: 2271.4285 -> 2311.111
$ws.Cells.Item(99, 10).Value = 2279.3333  # CRP!J99: 2412 -> 2279.3333
$ws.Cells.Item(99, 11).Value = 2311.111  # CRP!K99: 2271.4285 -> 2311.111
$ws.Cells.Item(99, 12).Value = 2279.3333  # CRP!L99: 2412 -> 2279.3333
$ws.Cells.Item(99, 13).Value = -813.1109999999999  # CRP!M99: -773.4285 -> -813.1109999999999
$ws.Cells.Item(99, 14).Value = -5275.3333  # CRP!N99: -5408 -> -5275.3333
$ws.Cells.Item(124, 8).Value = 26163  # CRP!H124: 0 -> 26163
$ws.Cells.Item(124, 10).Value = 26163  # CRP!J124: 0 -> 26163
$ws.Cells.Item(124, 12).Value = 26163  # CRP!L124: 0 -> 26163
$ws.Cells.Item(124, 14).Value = -31073  # CRP!N124: None -> -31073
$ws.Cells.Item(126, 8).Value = 2303.1667  # CRP!H126: 2302.6667 -> 2303.1667
$ws.Cells.Item(126, 9).Value = 2311.111  # CRP!I126: 2271.4285 -> 2311.111
$ws.Cells.Item(126, 10).Value = 2279.3333  # CRP!J126: 2412 -> 2279.3333
$ws.Cells.Item(126, 11).Value = 6933.333  # CRP!K126: 6814.2855 -> 6933.333
$ws.Cells.Item(126, 12).Value = 6837.999899999999  # CRP!L126: 7236 -> 6837.999899999999
$ws.Cells.Item(126, 13).Value = -4463.333  # CRP!M126: -4344.2855 -> -4463.333
$ws.Cells.Item(126, 14).Value = -11777.9999  # CRP!N126: -12176 -> -11777.9999
$ws.Cells.Item(132, 8).Value = 2500.6538  # CRP!H132: 2768.0435 -> 2500.6538
$ws.Cells.Item(132, 9).Value = 2226.5789  # CRP!I132: 2432.5293 -> 2226.5789
$ws.Cells.Item(132, 10).Value = 3244.5715  # CRP!J132: 3718.6667 -> 3244.5715
$ws.Cells.Item(132, 11).Value = 6679.736699999999  # CRP!K132: 7297.5879 -> 6679.736699999999
$ws.Cells.Item(132, 12).Value = 9733.7145  # CRP!L132: 11156.0001 -> 9733.7145
$ws.Cells.Item(132, 13).Value = -4149.736699999999  # CRP!M132: -4767.5879 -> -4149.736699999999
$ws.Cells.Item(132, 14).Value = -14793.7145  # CRP!N132: -16216.0001 -> -14793.7145
$ws.Cells.Item(134, 8).Value = 3661.4902  # CRP!H134: 3421.2727 -> 3661.4902
$ws.Cells.Item(134, 9).Value = 2500.4546  # CRP!I134: 2170.923 -> 2500.4546
$ws.Cells.Item(134, 11).Value = 7501.3638  # CRP!K134: 6512.768999999999 -> 7501.3638
$ws.Cells.Item(134, 13).Value = -4966.3638  # CRP!M134: -3977.768999999999 -> -4966.3638
$ws.Cells.Item(136, 8).Value = 1685500.9  # CRP!H136: 1654868.9 -> 1685500.9
$ws.Cells.Item(136, 9).Value = 2021626.8  # CRP!I136: 2067574 -> 2021626.8
$ws.Cells.Item(136, 10).Value = 4871.1113  # CRP!J136: 4048.5454 -> 4871.1113
$ws.Cells.Item(136, 11).Value = 6064880.4  # CRP!K136: 6202722 -> 6064880.4
$ws.Cells.Item(136, 12).Value = 14613.3339  # CRP!L136: 12145.6362 -> 14613.3339
$ws.Cells.Item(136, 13).Value = -6062330.4  # CRP!M136: -6200172 -> -6062330.4
$ws.Cells.Item(136, 14).Value = -19713.3339  # CRP!N136: -17245.6362 -> -19713.3339
$ws.Cells.Item(140, 8).Value = 35192.383  # CRP!H140: 43959.5 -> 35192.383
$ws.Cells.Item(140, 9).Value = 20709  # CRP!I140: 40709 -> 20709
$ws.Cells.Item(140, 10).Value = 36399.332  # CRP!J140: 44320.668 -> 36399.332
$ws.Cells.Item(140, 11).Value = 20709  # CRP!K140: 40709 -> 20709
$ws.Cells.Item(140, 12).Value = 36399.332  # CRP!L140: 44320.668 -> 36399.332
$ws.Cells.Item(140, 13).Value = -15529  # CRP!M140: -35529 -> -15529
$ws.Cells.Item(140, 14).Value = -46759.332  # CRP!N140: -54680.668 -> -46759.332

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 794.1087  # CUL!H113: 794.44446 -> 794.1087
$ws.Cells.Item(113, 9).Value = 792.0294  # CUL!I113: 792.42426 -> 792.0294
$ws.Cells.Item(113, 11).Value = 2376.0882  # CUL!K113: 2377.27278 -> 2376.0882
$ws.Cells.Item(113, 13).Value = -206.0882000000001  # CUL!M113: -207.2727800000002 -> -206.0882000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(92, 8).Value = 0  # GSM!H92: 14625.5 -> 0
$ws.Cells.Item(92, 10).Value = 0  # GSM!J92: 14625.5 -> 0
$ws.Cells.Item(92, 12).Value = 0  # GSM!L92: 14625.5 -> 0
$ws.Cells.Item(92, 14).ClearContents()  # GSM!N92: -18369.5 -> (removed)
$ws.Cells.Item(102, 8).Value = 6210.6665  # GSM!H102: 6340.8237 -> 6210.6665
$ws.Cells.Item(102, 9).Value = 5733.1665  # GSM!I102: 5890.909 -> 5733.1665
$ws.Cells.Item(102, 11).Value = 5733.1665  # GSM!K102: 5890.909 -> 5733.1665
$ws.Cells.Item(102, 13).Value = -4111.1665  # GSM!M102: -4268.909 -> -4111.1665
$ws.Cells.Item(126, 8).Value = 2800.84  # GSM!H126: 2755.16 -> 2800.84
$ws.Cells.Item(126, 9).Value = 1901.2  # GSM!I126: 1776.0834 -> 1901.2
$ws.Cells.Item(126, 10).Value = 3400.6  # GSM!J126: 3658.923 -> 3400.6
$ws.Cells.Item(126, 11).Value = 5703.6  # GSM!K126: 5328.2502 -> 5703.6
$ws.Cells.Item(126, 12).Value = 10201.8  # GSM!L126: 10976.769 -> 10201.8
$ws.Cells.Item(126, 13).Value = -3233.6  # GSM!M126: -2858.2502 -> -3233.6
$ws.Cells.Item(126, 14).Value = -15141.8  # GSM!N126: -15916.769 -> -15141.8
$ws.Cells.Item(132, 8).Value = 5534.409  # GSM!H132: 4525.891 -> 5534.409
$ws.Cells.Item(132, 9).Value = 4425.3335  # GSM!I132: 2471 -> 4425.3335
$ws.Cells.Item(132, 10).Value = 5950.3125  # GSM!J132: 6116.7744 -> 5950.3125
$ws.Cells.Item(132, 11).Value = 13276.0005  # GSM!K132: 7413 -> 13276.0005
$ws.Cells.Item(132, 12).Value = 17850.9375  # GSM!L132: 18350.3232 -> 17850.9375
$ws.Cells.Item(132, 13).Value = -10746.0005  # GSM!M132: -4883 -> -10746.0005
$ws.Cells.Item(132, 14).Value = -22910.9375  # GSM!N132: -23410.3232 -> -22910.9375

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 677.619  # LTW!H16: 734 -> 677.619
$ws.Cells.Item(16, 9).Value = 568.1111  # LTW!I16: 636.13336 -> 568.1111
$ws.Cells.Item(16, 10).Value = 1334.6666  # LTW!J16: 1101 -> 1334.6666
$ws.Cells.Item(16, 11).Value = 568.1111  # LTW!K16: 636.13336 -> 568.1111
$ws.Cells.Item(16, 12).Value = 1334.6666  # LTW!L16: 1101 -> 1334.6666
$ws.Cells.Item(16, 13).Value = -398.1111  # LTW!M16: -466.13336 -> -398.1111
$ws.Cells.Item(16, 14).Value = -1674.6666  # LTW!N16: -1441 -> -1674.6666
$ws.Cells.Item(122, 8).Value = 5987.0356  # LTW!H122: 7220.8096 -> 5987.0356
$ws.Cells.Item(122, 9).Value = 4542.316  # LTW!I122: 5450.2856 -> 4542.316
$ws.Cells.Item(122, 10).Value = 9037  # LTW!J122: 10761.857 -> 9037
$ws.Cells.Item(122, 11).Value = 13626.948  # LTW!K122: 16350.8568 -> 13626.948
$ws.Cells.Item(122, 12).Value = 27111  # LTW!L122: 32285.571 -> 27111
$ws.Cells.Item(122, 13).Value = -11176.948  # LTW!M122: -13900.8568 -> -11176.948
$ws.Cells.Item(122, 14).Value = -32011  # LTW!N122: -37185.571 -> -32011

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(92, 8).Value = 49550  # WVR!H92: 45050 -> 49550
$ws.Cells.Item(92, 10).Value = 49550  # WVR!J92: 45050 -> 49550
$ws.Cells.Item(92, 12).Value = 49550  # WVR!L92: 45050 -> 49550
$ws.Cells.Item(92, 14).Value = -54542  # WVR!N92: -50042 -> -54542
$ws.Cells.Item(122, 8).Value = 3182.577  # WVR!H122: 5064.9287 -> 3182.577
$ws.Cells.Item(122, 9).Value = 1242.5333  # WVR!I122: 1525 -> 1242.5333
$ws.Cells.Item(122, 10).Value = 5828.091  # WVR!J122: 9784.833000000001 -> 5828.091
$ws.Cells.Item(122, 11).Value = 3727.5999  # WVR!K122: 4575 -> 3727.5999
$ws.Cells.Item(122, 12).Value = 17484.273  # WVR!L122: 29354.499 -> 17484.273
$ws.Cells.Item(122, 13).Value = -1277.5999  # WVR!M122: -2125 -> -1277.5999
$ws.Cells.Item(122, 14).Value = -22384.273  # WVR!N122: -34254.499 -> -22384.273
$ws.Cells.Item(125, 8).Value = 31333  # WVR!H125: 31081 -> 31333
$ws.Cells.Item(125, 10).Value = 31333  # WVR!J125: 31081 -> 31333
$ws.Cells.Item(125, 12).Value = 31333  # WVR!L125: 31081 -> 31333
$ws.Cells.Item(125, 14).Value = -41173  # WVR!N125: -40921 -> -41173
